# Fixing new OWID attributes and Micronesia
#
# 1. Insert a new row for "Micronesia" at row 132 (pushing Moldova and
#    everything below it down by one row).
# 2. Update the "alpha2" defined name so its range keeps up with the
#    newly added row.
# 3. Move the sheet's visible window / selection to where the edit
#    happened.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Micronesia" row ------------------------------------
$ws.Rows.Item(132).Insert()

$ws.Range("A132").Value = "Micronesia"
$ws.Range("B132").Value = "FM"
$ws.Range("C132").Value = "FSM"
$ws.Range("D132").Value = "Oceania"
$ws.Range("E132").Value = 112640

# --- 2. Keep the "alpha2" named range in sync with the extra row -----------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "Tabelle1!alpha2") {
        $n.RefersTo = "=Tabelle1!`$A`$1:`$C`$219"
    }
}

# --- 3. Move the window / selection to reflect where we were working -------
$win = $excel.ActiveWindow
$win.ScrollRow = 114
$win.ScrollColumn = 1
$ws.Range("A132:E132").Select()
